$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = [double]"7.16e-19"
$ws.Cells.Item(2, 4).Value = [double]"1.352"
$ws.Cells.Item(2, 5).Value = [double]"0"
$ws.Cells.Item(2, 6).Value = [double]"2"
$ws.Cells.Item(3, 3).Value = [double]"0.39"
$ws.Cells.Item(3, 4).Value = [double]"1.087"
$ws.Cells.Item(3, 5).Value = [double]"8.346"
$ws.Cells.Item(3, 6).Value = [double]"14"
$ws.Cells.Item(4, 3).Value = [double]"0.403"
$ws.Cells.Item(4, 4).Value = [double]"1.087"
$ws.Cells.Item(4, 5).Value = [double]"3.376"
$ws.Cells.Item(4, 6).Value = [double]"14"
$ws.Cells.Item(5, 3).Value = [double]"0.4956"
$ws.Cells.Item(5, 4).Value = [double]"1.283"
$ws.Cells.Item(5, 5).Value = [double]"0.4391"
$ws.Cells.Item(5, 6).Value = [double]"17"
$ws.Cells.Item(6, 3).Value = [double]"0.3475"
$ws.Cells.Item(6, 4).Value = [double]"1.216"
$ws.Cells.Item(6, 5).Value = [double]"1.154"
$ws.Cells.Item(6, 6).Value = [double]"14"
$ws.Cells.Item(7, 3).Value = [double]"0.01168"
$ws.Cells.Item(7, 4).Value = [double]"1.925"
$ws.Cells.Item(7, 5).Value = [double]"2.734"
$ws.Cells.Item(7, 6).Value = [double]"7"
$ws.Cells.Item(8, 3).Value = [double]"0.7325"
$ws.Cells.Item(8, 4).Value = [double]"1.411"
$ws.Cells.Item(8, 5).Value = [double]"4.85"
$ws.Cells.Item(8, 6).Value = [double]"21"
$ws.Cells.Item(9, 3).Value = [double]"0.07799"
$ws.Cells.Item(9, 4).Value = [double]"2.249"
$ws.Cells.Item(9, 6).Value = [double]"6"
$ws.Cells.Item(10, 3).Value = [double]"0.165"
$ws.Cells.Item(10, 4).Value = [double]"1.887"
$ws.Cells.Item(10, 6).Value = [double]"10"
$ws.Cells.Item(11, 3).Value = [double]"0.158"
$ws.Cells.Item(11, 4).Value = [double]"1.893"
$ws.Cells.Item(11, 6).Value = [double]"11"
$ws.Cells.Item(12, 3).Value = [double]"0.2192"
$ws.Cells.Item(12, 4).Value = [double]"1.961"
$ws.Cells.Item(12, 5).Value = [double]"5.746"
$ws.Cells.Item(12, 6).Value = [double]"13"
$ws.Cells.Item(13, 3).Value = [double]"0.244"
$ws.Cells.Item(13, 4).Value = [double]"2.067"
$ws.Cells.Item(13, 5).Value = [double]"0.1693"
$ws.Cells.Item(13, 6).Value = [double]"14"
$ws.Cells.Item(14, 3).Value = [double]"0.2771"
$ws.Cells.Item(14, 4).Value = [double]"2.147"
$ws.Cells.Item(14, 5).Value = [double]"1.517"
$ws.Cells.Item(14, 6).Value = [double]"16"
$ws.Cells.Item(15, 3).Value = [double]"0.8951"
$ws.Cells.Item(15, 4).Value = [double]"1.877"
$ws.Cells.Item(15, 5).Value = [double]"2.119"
$ws.Cells.Item(15, 6).Value = [double]"21"
$ws.Cells.Item(16, 3).Value = [double]"0.00134"
$ws.Cells.Item(16, 4).Value = [double]"2.986"
$ws.Cells.Item(16, 6).Value = [double]"6"
$ws.Cells.Item(17, 3).Value = [double]"8.933000000000001e-19"
$ws.Cells.Item(17, 4).Value = [double]"4.673"
$ws.Cells.Item(17, 5).Value = [double]"0"
$ws.Cells.Item(17, 6).Value = [double]"2"
$ws.Cells.Item(18, 3).Value = [double]"0.003378"
$ws.Cells.Item(18, 4).Value = [double]"3.301"
$ws.Cells.Item(18, 6).Value = [double]"8"
$ws.Cells.Item(19, 3).Value = [double]"0.0537"
$ws.Cells.Item(19, 4).Value = [double]"2.855"
$ws.Cells.Item(19, 5).Value = [double]"1.833"
$ws.Cells.Item(19, 6).Value = [double]"12"
$ws.Cells.Item(20, 3).Value = [double]"0.03501"
$ws.Cells.Item(20, 4).Value = [double]"3.166"
$ws.Cells.Item(20, 5).Value = [double]"3.975"
$ws.Cells.Item(20, 6).Value = [double]"10"
$ws.Cells.Item(21, 3).Value = [double]"0.03476"
$ws.Cells.Item(21, 4).Value = [double]"3.315"
$ws.Cells.Item(21, 5).Value = [double]"2.251"
$ws.Cells.Item(21, 6).Value = [double]"11"
$ws.Cells.Item(22, 3).Value = [double]"0.03365"
$ws.Cells.Item(22, 4).Value = [double]"3.471"
$ws.Cells.Item(22, 5).Value = [double]"2.27"
$ws.Cells.Item(22, 6).Value = [double]"12"
$ws.Cells.Item(23, 3).Value = [double]"0.07542"
$ws.Cells.Item(23, 4).Value = [double]"3.46"
$ws.Cells.Item(23, 5).Value = [double]"2.956"
$ws.Cells.Item(23, 6).Value = [double]"16"
$ws.Cells.Item(24, 3).Value = [double]"0.1338"
$ws.Cells.Item(24, 4).Value = [double]"3.465"
$ws.Cells.Item(24, 5).Value = [double]"1.083"
$ws.Cells.Item(24, 6).Value = [double]"18"
$ws.Cells.Item(25, 3).Value = [double]"0.1851"
$ws.Cells.Item(25, 4).Value = [double]"3.375"
$ws.Cells.Item(25, 5).Value = [double]"0.6208"
$ws.Cells.Item(25, 6).Value = [double]"21"
$ws.Cells.Item(26, 3).Value = [double]"0.429"
$ws.Cells.Item(26, 4).Value = [double]"7.131"
$ws.Cells.Item(26, 5).Value = [double]"0"
$ws.Cells.Item(26, 6).Value = [double]"1"
$ws.Cells.Item(27, 3).Value = [double]"1.262e-29"
$ws.Cells.Item(27, 4).Value = [double]"6.608"
$ws.Cells.Item(27, 5).Value = [double]"0"
$ws.Cells.Item(27, 6).Value = [double]"4"
$ws.Cells.Item(28, 3).Value = [double]"0.001944"
$ws.Cells.Item(28, 4).Value = [double]"5.613"
$ws.Cells.Item(28, 5).Value = [double]"1.064"
$ws.Cells.Item(28, 6).Value = [double]"6"
$ws.Cells.Item(29, 3).Value = [double]"0.002924"
$ws.Cells.Item(29, 4).Value = [double]"5.429"
$ws.Cells.Item(29, 6).Value = [double]"6"
$ws.Cells.Item(30, 3).Value = [double]"0.2944"
$ws.Cells.Item(30, 4).Value = [double]"8.938000000000001"
$ws.Cells.Item(30, 5).Value = [double]"0"
$ws.Cells.Item(30, 6).Value = [double]"1"
$ws.Cells.Item(31, 3).Value = [double]"0.01533"
$ws.Cells.Item(31, 4).Value = [double]"5.019"
$ws.Cells.Item(31, 5).Value = [double]"0.9826"
$ws.Cells.Item(31, 6).Value = [double]"13"
$ws.Cells.Item(32, 3).Value = [double]"0.1035"
$ws.Cells.Item(32, 4).Value = [double]"4.556"
$ws.Cells.Item(32, 5).Value = [double]"0.764"
$ws.Cells.Item(32, 6).Value = [double]"18"
$ws.Cells.Item(33, 3).Value = [double]"0.0824"
$ws.Cells.Item(33, 4).Value = [double]"4.61"
$ws.Cells.Item(33, 5).Value = [double]"1.01"
$ws.Cells.Item(33, 6).Value = [double]"18"
$ws.Cells.Item(34, 3).Value = [double]"0.1264"
$ws.Cells.Item(34, 4).Value = [double]"4.457"
$ws.Cells.Item(34, 5).Value = [double]"0.301"
$ws.Cells.Item(34, 6).Value = [double]"21"
$ws.Cells.Item(35, 3).Value = [double]"0.104"
$ws.Cells.Item(35, 4).Value = [double]"4.554"
$ws.Cells.Item(35, 5).Value = [double]"0.6122"
$ws.Cells.Item(35, 6).Value = [double]"21"
$ws.Cells.Item(36, 3).Value = [double]"2.312e-19"
$ws.Cells.Item(36, 4).Value = [double]"7.081"
$ws.Cells.Item(36, 5).Value = [double]"0"
$ws.Cells.Item(36, 6).Value = [double]"2"
$ws.Cells.Item(37, 3).Value = [double]"0.01241"
$ws.Cells.Item(37, 4).Value = [double]"5.121"
$ws.Cells.Item(37, 6).Value = [double]"10"
$ws.Cells.Item(38, 3).Value = [double]"0.04219"
$ws.Cells.Item(38, 4).Value = [double]"5.038"
$ws.Cells.Item(38, 5).Value = [double]"0.8683999999999999"
$ws.Cells.Item(38, 6).Value = [double]"14"
$ws.Cells.Item(39, 3).Value = [double]"0.03904"
$ws.Cells.Item(39, 4).Value = [double]"5.085"
$ws.Cells.Item(39, 6).Value = [double]"15"
$ws.Cells.Item(40, 3).Value = [double]"0.06263000000000001"
$ws.Cells.Item(40, 4).Value = [double]"5.082"
$ws.Cells.Item(40, 5).Value = [double]"1.165"
$ws.Cells.Item(40, 6).Value = [double]"17"
$ws.Cells.Item(41, 3).Value = [double]"0.1733"
$ws.Cells.Item(41, 4).Value = [double]"4.953"
$ws.Cells.Item(41, 5).Value = [double]"0.07194"
$ws.Cells.Item(41, 6).Value = [double]"21"
$ws.Cells.Item(42, 3).Value = [double]"0.00175"
$ws.Cells.Item(42, 4).Value = [double]"5.667"
$ws.Cells.Item(42, 6).Value = [double]"6"
$ws.Cells.Item(43, 3).Value = [double]"0.02719"
$ws.Cells.Item(43, 4).Value = [double]"5.177"
$ws.Cells.Item(43, 6).Value = [double]"11"
$ws.Cells.Item(44, 3).Value = [double]"0.6196"
$ws.Cells.Item(44, 4).Value = [double]"5.13"
$ws.Cells.Item(44, 5).Value = [double]"5.344"
$ws.Cells.Item(44, 6).Value = [double]"20"
$ws.Cells.Item(45, 3).Value = [double]"0.5965"
$ws.Cells.Item(45, 4).Value = [double]"5.141"
$ws.Cells.Item(45, 5).Value = [double]"2.722"
$ws.Cells.Item(45, 6).Value = [double]"21"
$ws.Cells.Item(46, 3).Value = [double]"0.00479"
$ws.Cells.Item(46, 4).Value = [double]"6.554"
$ws.Cells.Item(46, 6).Value = [double]"5"
$ws.Cells.Item(47, 3).Value = [double]"0.06633"
$ws.Cells.Item(47, 4).Value = [double]"5.719"
$ws.Cells.Item(47, 5).Value = [double]"1.158"
$ws.Cells.Item(47, 6).Value = [double]"9"
$ws.Cells.Item(48, 3).Value = [double]"1.252e-19"
$ws.Cells.Item(48, 4).Value = [double]"7.863"
$ws.Cells.Item(48, 5).Value = [double]"0"
$ws.Cells.Item(48, 6).Value = [double]"3"
$ws.Cells.Item(49, 3).Value = [double]"0.1792"
$ws.Cells.Item(49, 4).Value = [double]"5.803"
$ws.Cells.Item(49, 5).Value = [double]"0.9627"
$ws.Cells.Item(49, 6).Value = [double]"13"
$ws.Cells.Item(50, 3).Value = [double]"1.037"
$ws.Cells.Item(50, 4).Value = [double]"5.371"
$ws.Cells.Item(50, 6).Value = [double]"21"
$ws.Cells.Item(51, 3).Value = [double]"1.367"
$ws.Cells.Item(51, 4).Value = [double]"9.071"
$ws.Cells.Item(51, 5).Value = [double]"0"
$ws.Cells.Item(51, 6).Value = [double]"1"
$ws.Cells.Item(52, 3).Value = [double]"1.349e-18"
$ws.Cells.Item(52, 4).Value = [double]"7.858"
$ws.Cells.Item(52, 5).Value = [double]"0"
$ws.Cells.Item(52, 6).Value = [double]"3"
$ws.Cells.Item(53, 3).Value = [double]"0.03151"
$ws.Cells.Item(53, 4).Value = [double]"5.921"
$ws.Cells.Item(53, 5).Value = [double]"0.6631"
$ws.Cells.Item(53, 6).Value = [double]"8"
$ws.Cells.Item(54, 3).Value = [double]"0.08048"
$ws.Cells.Item(54, 4).Value = [double]"5.491"
$ws.Cells.Item(54, 6).Value = [double]"11"
$ws.Cells.Item(55, 3).Value = [double]"5.889e-19"
$ws.Cells.Item(55, 4).Value = [double]"6.233"
$ws.Cells.Item(55, 5).Value = [double]"0"
$ws.Cells.Item(55, 6).Value = [double]"3"
$ws.Cells.Item(56, 3).Value = [double]"0.002066"
$ws.Cells.Item(56, 4).Value = [double]"5.591"
$ws.Cells.Item(56, 6).Value = [double]"5"
$ws.Cells.Item(57, 3).Value = [double]"0.3455"
$ws.Cells.Item(57, 4).Value = [double]"5.341"
$ws.Cells.Item(57, 5).Value = [double]"0.2023"
$ws.Cells.Item(57, 6).Value = [double]"17"
$ws.Cells.Item(58, 3).Value = [double]"0.4082"
$ws.Cells.Item(58, 4).Value = [double]"5.227"
$ws.Cells.Item(58, 5).Value = [double]"0.8213"
$ws.Cells.Item(58, 6).Value = [double]"18"
$ws.Cells.Item(59, 3).Value = [double]"0.4576"
$ws.Cells.Item(59, 4).Value = [double]"5.157"
$ws.Cells.Item(59, 5).Value = [double]"0.8869"
$ws.Cells.Item(60, 3).Value = [double]"5.74e-05"
$ws.Cells.Item(60, 4).Value = [double]"5.142"
$ws.Cells.Item(60, 6).Value = [double]"5"
$ws.Cells.Item(61, 3).Value = [double]"0.01524"
$ws.Cells.Item(61, 4).Value = [double]"4.605"
$ws.Cells.Item(61, 6).Value = [double]"7"
$ws.Cells.Item(62, 3).Value = [double]"0.03197"
$ws.Cells.Item(62, 4).Value = [double]"4.288"
$ws.Cells.Item(62, 6).Value = [double]"11"
$ws.Cells.Item(63, 3).Value = [double]"0.01792"
$ws.Cells.Item(63, 4).Value = [double]"4.276"
$ws.Cells.Item(63, 5).Value = [double]"2.296"
$ws.Cells.Item(63, 6).Value = [double]"11"
$ws.Cells.Item(64, 3).Value = [double]"8.577e-19"
$ws.Cells.Item(64, 4).Value = [double]"6.08"
$ws.Cells.Item(64, 5).Value = [double]"0"
$ws.Cells.Item(64, 6).Value = [double]"3"
$ws.Cells.Item(65, 3).Value = [double]"0.7052"
$ws.Cells.Item(65, 4).Value = [double]"10.28"
$ws.Cells.Item(65, 5).Value = [double]"0"
$ws.Cells.Item(65, 6).Value = [double]"1"
$ws.Cells.Item(66, 3).Value = [double]"0.0009496"
$ws.Cells.Item(66, 4).Value = [double]"6.992"
$ws.Cells.Item(66, 6).Value = [double]"5"
$ws.Cells.Item(67, 3).Value = [double]"0.01128"
$ws.Cells.Item(67, 4).Value = [double]"6.113"
$ws.Cells.Item(67, 6).Value = [double]"6"
$ws.Cells.Item(68, 3).Value = [double]"0.01327"
$ws.Cells.Item(68, 4).Value = [double]"5.244"
$ws.Cells.Item(68, 5).Value = [double]"3.092"
$ws.Cells.Item(68, 6).Value = [double]"9"
$ws.Cells.Item(69, 3).Value = [double]"0.01674"
$ws.Cells.Item(69, 4).Value = [double]"5.044"
$ws.Cells.Item(69, 5).Value = [double]"0.1013"
$ws.Cells.Item(69, 6).Value = [double]"12"
$ws.Cells.Item(70, 3).Value = [double]"0.0002169"
$ws.Cells.Item(70, 4).Value = [double]"6.986"
$ws.Cells.Item(70, 6).Value = [double]"5"
$ws.Cells.Item(71, 3).Value = [double]"0.08319"
$ws.Cells.Item(71, 4).Value = [double]"5.386"
$ws.Cells.Item(71, 5).Value = [double]"0.6427"
$ws.Cells.Item(71, 6).Value = [double]"16"
$ws.Cells.Item(72, 3).Value = [double]"0.116"
$ws.Cells.Item(72, 4).Value = [double]"5.282"
$ws.Cells.Item(72, 5).Value = [double]"0.2679"
$ws.Cells.Item(72, 6).Value = [double]"18"
